$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.738864421844482
$ws.Range("B1").Value = 1.843623876571655
$ws.Range("C1").Value = 1.990047097206116
$ws.Range("D1").Value = 2.850121974945068
$ws.Range("E1").Value = 3.397604942321777
